$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '23.463.45'
$ws.Range("E2").Value = '  +1.98%  '
$ws.Range("D3").Value = '1.630.94'
$ws.Range("E3").Value = '  +2.95%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9977'
$ws.Range("E4").Value = '  -0.56%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '307.92'
$ws.Range("E5").Value = '  +2.51%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9978'
$ws.Range("E6").Value = '  -0.51%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3785'
$ws.Range("E7").Value = '  +0.86%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '53.03'
$ws.Range("E8").Value = '  +5.07%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3666'
$ws.Range("E9").Value = '  +2.70%  '
$ws.Range("E10").Value = '  +6.27%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.08207'
$ws.Range("E11").Value = '  +2.97%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.9978'
$ws.Range("E12").Value = '  -0.56%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '23.35'
$ws.Range("E13").Value = '  +7.40%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.687'
$ws.Range("E14").Value = '  +3.88%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.00001272'
$ws.Range("E15").Value = '  +4.56%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '7.470'
$ws.Range("E16").Value = '  +2.79%  '
$ws.Range("D17").Value = '1.625.15'
$ws.Range("E17").Value = '  +2.13%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '94.95'
$ws.Range("E18").Value = '  +3.10%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06948'
$ws.Range("E19").Value = '  +2.85%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '18.43'
$ws.Range("E20").Value = '  +4.02%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.603'
$ws.Range("E21").Value = '  +3.38%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.9981'
$ws.Range("E22").Value = '  -0.42%  '
$ws.Range("B23").Value = 'Cosmos'
$ws.Range("C23").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '13.02'
$ws.Range("E23").Value = '  +2.38%  '
$ws.Range("B24").Value = 'WrappedBTC'
$ws.Range("C24").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D24").Value = '23.471.55'
$ws.Range("E24").Value = '  +2.02%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.133'
$ws.Range("E25").Value = '  +14.05%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.435'
$ws.Range("E26").Value = '  +2.56%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '21.42'
$ws.Range("E27").Value = '  +3.91%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '151.32'
$ws.Range("E28").Value = '  +3.21%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.315'
$ws.Range("E29").Value = '  +2.28%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '136.59'
$ws.Range("E30").Value = '  +3.71%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.433'
$ws.Range("E31").Value = '  +5.33%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.928'
$ws.Range("E32").Value = '  +6.76%  '
$ws.Range("D33").Value = '1.807.66'
$ws.Range("E33").Value = '  +2.47%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.9810'
$ws.Range("E34").Value = '  +5.62%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.02811'
$ws.Range("E35").Value = '  +5.81%  '
$ws.Range("E36").Value = '  +5.87%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.07475'
$ws.Range("E37").Value = '  +2.36%  '
$ws.Range("E38").Value = '  +4.74%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.2538'
$ws.Range("E39").Value = '  +2.45%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.08851'
$ws.Range("E40").Value = '  +1.39%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.405'
$ws.Range("E41").Value = '  +5.67%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.7182'
$ws.Range("E42").Value = '  +5.32%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '12.85'
$ws.Range("E43").Value = '  +9.33%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '16.30'
$ws.Range("E44").Value = '  +11.09%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.6634'
$ws.Range("E45").Value = '  +4.72%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.370'
$ws.Range("E46").Value = '  +6.10%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.036'
$ws.Range("E47").Value = '  +1.49%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.9970'
$ws.Range("E48").Value = '  -0.49%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.08029'
$ws.Range("E49").Value = '  +2.02%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '132.74'
$ws.Range("E50").Value = '  +1.47%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.220'
$ws.Range("E51").Value = '  +3.41%  '
